# Apply "repull data" update to the dSF column (F) values.
# These values come from a re-pull of upstream data; they are written directly
# as literal numbers (not formulas) to match the source workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dSF = @{
    2  = -2
    3  = 1
    5  = -1
    6  = 2
    7  = 1
    8  = 6
    9  = -3
    10 = -1
    11 = 1
    12 = -4
    13 = -2
    14 = 4
    15 = -1
    16 = 8
    17 = 9
    19 = -2
    20 = 2
    22 = -2
    23 = -1
    24 = -1
    25 = -3
    26 = -2
    28 = -2
    29 = -1
    30 = -1
}

foreach ($row in $dSF.Keys) {
    $ws.Range("F$row").Value = $dSF[$row]
}
